$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = @{
    2  = "sports_club_boys"
    3  = "sports_uil_boys"
    4  = "sports_club_girls"
    5  = "sports_uil_girls"
    6  = "sports_uil_boys"
    7  = "sports_uil_girls"
    8  = "sports_uil_boys"
    9  = "sports_uil_girls"
    10 = "sports_uil_boys"
    11 = "sports_uil_girls"
    12 = "sports_uil_boys"
    13 = "sports_uil_girls"
    14 = "sports_uil_boys"
    15 = "sports_uil_girls"
    16 = "sports_uil_coed"
    17 = "sports_uil_boys"
    18 = "sports_uil_girls"
}

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 2).Value = $newB[$row]
    $ws.Cells.Item($row, 3).Value = "Water Polo"
}
